$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (result of sharedStrings table reordering in the source diff) ---
# Indonesia overtakes Arabia Saudita in case count ranking
$ws.Range("A40").Value = "Indonesia"
$ws.Range("A41").Value = "Arabia Saudita"

# Azerbaiyan overtakes Principado de Andorra and Costa Rica
$ws.Range("A79").Value = "Azerbaiyan"
$ws.Range("A80").Value = "Principado de Andorra"
$ws.Range("A81").Value = "Costa Rica"

# Camerun overtakes Afganistan and San Marino
$ws.Range("A91").Value = "Camerun"
$ws.Range("A92").Value = "Afganistan"
$ws.Range("A93").Value = "San Marino"

# --- Updated "last refreshed" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 10:50"

# --- Refreshed COVID-19 statistics (columns B:H) for affected rows ---
# Row 24
$ws.Range("B24").Value = 4898
$ws.Range("C24").Value = 21
$ws.Range("E24").Value = 4840
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 45
# Row 40
$ws.Range("B40").Value = 1790
$ws.Range("C40").Value = 113
$ws.Range("D40").Value = 112
$ws.Range("E40").Value = 1508
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 13
$ws.Range("H40").Value = 170
# Row 41
$ws.Range("B41").Value = 1720
$ws.Range("D41").Value = 264
$ws.Range("E41").Value = 1440
$ws.Range("F41").Value = 31
$ws.Range("H41").Value = 16
# Row 71
$ws.Range("B71").Value = 512
$ws.Range("C71").Value = 53
$ws.Range("D71").Value = 20
$ws.Range("E71").Value = 477
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 15
# Row 72
$ws.Range("E72").Value = 421
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 15
# Row 76
$ws.Range("F76").Value = 65
# Row 79
$ws.Range("B79").Value = 400
$ws.Range("C79").Value = 41
$ws.Range("D79").Value = 26
$ws.Range("E79").Value = 369
$ws.Range("F79").Value = 7
$ws.Range("H79").Value = 5
# Row 80
$ws.Range("B80").Value = 390
$ws.Range("D80").Value = 10
$ws.Range("E80").Value = 366
$ws.Range("F80").Value = 12
$ws.Range("H80").Value = 14
# Row 81
$ws.Range("B81").Value = 375
$ws.Range("D81").Value = 4
$ws.Range("E81").Value = 369
$ws.Range("F81").Value = 9
$ws.Range("H81").Value = 2
# Row 91
$ws.Range("B91").Value = 255
$ws.Range("C91").Value = 22
$ws.Range("D91").Value = 10
$ws.Range("E91").Value = 239
$ws.Range("H91").Value = 6
# Row 92
$ws.Range("B92").Value = 239
$ws.Range("C92").Value = 2
$ws.Range("D92").Value = 5
$ws.Range("E92").Value = 230
$ws.Range("F92").Value = 0
$ws.Range("H92").Value = 4
# Row 93
$ws.Range("B93").Value = 236
$ws.Range("D93").Value = 13
$ws.Range("E93").Value = 195
$ws.Range("F93").Value = 16
$ws.Range("H93").Value = 28
# Row 94
$ws.Range("D94").Value = 57
$ws.Range("E94").Value = 173
